$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: insert $content at absolute position $pos, then wrap the newly
# inserted span with a throw-away bookmark so that it stays a run distinct
# from its neighbours even while later edits happen elsewhere in the
# document. All such helper bookmarks are collected in $script:tempMarks and
# removed right at the end, once every textual edit has been made (deleting
# a bookmark afterwards does not cause Word to re-merge the runs it already
# separated).
# ---------------------------------------------------------------------------
$script:tempMarks = New-Object System.Collections.ArrayList
$script:markCounter = 0

function New-RunBoundaryMark([int]$pos) {
    # Zero-width bookmark: keeps the run boundary at $pos without inserting
    # any visible text (used when the text on both sides already exists,
    # e.g. separating "...key_" from "return...").
    $script:markCounter += 1
    $name = "zzKeepMark" + $script:markCounter
    $r = $d.Range($pos, $pos)
    $d.Bookmarks.Add($name, $r) | Out-Null
    [void]$script:tempMarks.Add($name)
}

function Insert-IsolatedText([int]$pos, [string]$content) {
    # Inserts $content at $pos and wraps it with a bookmark so it remains
    # its own run (distinct from the text before and after it).
    $ins = $d.Range($pos, $pos)
    $ins.InsertBefore($content)
    $script:markCounter += 1
    $name = "zzKeepMark" + $script:markCounter
    $r = $d.Range($pos, $pos + $content.Length)
    $d.Bookmarks.Add($name, $r) | Out-Null
    [void]$script:tempMarks.Add($name)
}

function Remove-TempMarks() {
    foreach ($n in $script:tempMarks) {
        if ($d.Bookmarks.Exists($n)) {
            $d.Bookmarks($n).Delete()
        }
    }
    $script:tempMarks.Clear()
}

# ---------------------------------------------------------------------------
# 1) Date paragraph: "04/13/2017" -> "04/24" / _GoBack bookmark / "/2017"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("04/13/2017", $true, $false, $false, $false, $false,
                         $true, 1, $false, "04/24/2017", 2) | Out-Null

$p1 = $d.Paragraphs(1).Range
$p1Start = $p1.Start
$p1Text = $p1.Text
$splitPos = $p1Start + $p1Text.IndexOf("/2017")
$d.Bookmarks.Add("_GoBack", $d.Range($splitPos, $splitPos)) | Out-Null

# ---------------------------------------------------------------------------
# 2) "Version 1.0" -> "Version 1.1"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Version 1.0", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Version 1.1", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) "Record each input as log file. ..." -> split "as" into its own run,
#    then replace its text with "in a".
# ---------------------------------------------------------------------------
$pRec = $d.Paragraphs(4).Range
$recStart = $pRec.Start
$recText = $pRec.Text
$asIdx = $recText.IndexOf("as log")
$asRange = $d.Range($recStart + $asIdx, $recStart + $asIdx + 2)
$asMarkName = "zzAsMark"
$d.Bookmarks.Add($asMarkName, $asRange) | Out-Null
$d.Bookmarks($asMarkName).Range.Text = "in a"
$d.Bookmarks($asMarkName).Delete()

# ---------------------------------------------------------------------------
# 4) "c:/qfl2017/checkin.txt," -> "c:/qfl" / "/" / "2017/checkin.txt,"
# ---------------------------------------------------------------------------
$pChk = $d.Paragraphs(5).Range
$chkStart = $pChk.Start
$chkText = $pChk.Text
$chkSlashPos = $chkStart + $chkText.IndexOf("2017/checkin")
Insert-IsolatedText $chkSlashPos "/"

# ---------------------------------------------------------------------------
# 5) "c:/qfl2017/key-return.txt," -> "c:/qfl" / "/" / "2017/key_" / "return.txt,"
#    and the trailing _GoBack bookmark that used to sit in this paragraph is
#    gone (it was relocated to the date paragraph in step 1).
# ---------------------------------------------------------------------------
$pKey = $d.Paragraphs(6).Range
$keyStart = $pKey.Start
$keyText = $pKey.Text

# 5a) Mark the boundary between the (future) underscore and "return" so the
#     dash->underscore replacement does not get merged with "return.txt,".
$returnIdx = $keyText.IndexOf("return.txt")
New-RunBoundaryMark ($keyStart + $returnIdx)

# 5b) Replace "-" with "_" (same length, stays put, protected on the right
#     by the boundary mark placed above).
$dashIdx = $keyText.IndexOf("-")
$dashRange = $d.Range($keyStart + $dashIdx, $keyStart + $dashIdx + 1)
$dashRange.Text = "_"

# 5c) Insert "/" between "qfl" and "2017" as its own isolated run.
$keyText2 = $d.Paragraphs(6).Range.Text
$keyStart2 = $d.Paragraphs(6).Range.Start
$keySlashPos = $keyStart2 + $keyText2.IndexOf("2017/key")
Insert-IsolatedText $keySlashPos "/"

# ---------------------------------------------------------------------------
# All textual edits that needed run isolation are done: drop every helper
# bookmark now (this will not re-merge the already-separated runs).
# ---------------------------------------------------------------------------
Remove-TempMarks

# ---------------------------------------------------------------------------
# 6) Insert a new list item after "In key-return, add attendee number..."
# ---------------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "In key-return, add attendee number*") {
        $target = $i
        break
    }
}
$srcPara = $d.Paragraphs($target).Range
$srcPara.InsertParagraphAfter()
$newPara = $d.Paragraphs($target + 1).Range
$newPara.Text = "Add control button in checkin to allow full payment update. "
